# ProjectCreation.xlsx edit — rename sheet, refresh credentials, add
# project-fields columns (Name/Description/Start date) on the "Mandatory
# name" sheet, resize the new columns, and update the active selection.

$wb = $excel.ActiveWorkbook

# 1) Rename the first sheet from "Existing" to "Mandatory name"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Mandatory name"

# 2) Refresh the username/password sample data on the renamed sheet
$ws1.Range("A2").Value = "angel_dlr"
$ws1.Range("B2").Value = 'Pa$$w0rd'
$ws1.Range("A3").Value = "angel_dlr"
$ws1.Range("B3").Value = 'Pa$$w0rd'

# 3) Add the project-creation header row and sample rows
$ws1.Range("C1").Value = "Name"
$ws1.Range("D1").Value = "Description"
$ws1.Range("E1").Value = "Start date"

$ws1.Range("C2").Value = "Testing Project"
$ws1.Range("D2").Value = "During this project we expect on having successful tests that will prove that our software works properly."
$ws1.Range("E2").Value = "18 October 2021"

$ws1.Range("D3").Value = "Create optimal test cases to perform for the correct operation of the program"
$ws1.Range("E3").Value = "15 October 2021"

# 4) Widen the new Name/Description columns
$ws1.Columns.Item(3).ColumnWidth = 21.8
$ws1.Columns.Item(4).ColumnWidth = 24.65

# 5) Update the active selection on the renamed sheet
$ws1.Range("G21").Select() | Out-Null
